# Atualização de bases das ligas, do dia: 19-04-2024 às 21:40
#
# For each of the listed row pairs, the data in columns B..AC (everything
# except the pandas-style index column A) was swapped between the two
# rows. In addition, the very last data row (158, a still-unplayed /
# unmatched fixture) was removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param(
        [int]$RowA,
        [int]$RowB
    )

    $rangeA = $ws.Range("B$RowA`:AC$RowA")
    $rangeB = $ws.Range("B$RowB`:AC$RowB")

    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2

    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

# Row pairs whose content got swapped (positional args — named params are
# not reliably bound by this PowerShell host)
Swap-RowData 9   10
Swap-RowData 76  77
Swap-RowData 87  88
Swap-RowData 111 112
Swap-RowData 122 123

# Remove the trailing row 158 entirely so the used range becomes A1:AC157.
$ws.Range("A158:AC158").EntireRow.Delete()
